$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.466.29'
$ws.Range('E2').Value = '  +2.78%  '
$ws.Range('D3').Value = '2.314.73'
$ws.Range('E3').Value = '  +1.90%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.46'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.82'
$ws.Range('E6').Value = '  +7.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.534'
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +8.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.85'
$ws.Range('E10').Value = '  +5.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.85'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0814'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.03'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = '2.672.48'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('D17').Value = '2.311.33'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.814'
$ws.Range('E18').Value = '  +3.35%  '
$ws.Range('D19').Value = '43.379.58'
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.23'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').Value = '0.0₃0930'
$ws.Range('E21').Value = '  +2.64%  '
$ws.Range('E22').Value = '  +3.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.34'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '242.70'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('E25').Value = '  +3.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.62'
$ws.Range('E26').Value = '  +1.10%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.83'
$ws.Range('E28').Value = '  +5.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.19'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.66'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.76'
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.30'
$ws.Range('E33').Value = '  +1.20%  '
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.34'
$ws.Range('E35').Value = '  +4.21%  '
$ws.Range('E36').Value = '  +6.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0746'
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.57'
$ws.Range('E39').Value = '  +9.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.89'
$ws.Range('E40').Value = '  +3.83%  '
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.116'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.73'
$ws.Range('E43').Value = '  +19.96%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0294'
$ws.Range('E44').Value = '  +3.70%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.997.46'
$ws.Range('E45').Value = '  +2.67%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.10'
$ws.Range('E46').Value = '  +1.86%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.12'
$ws.Range('E47').Value = '  +5.94%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.03'
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '57.01'
$ws.Range('E49').Value = '  +5.26%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.96'
$ws.Range('E50').Value = '  +2.17%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.60'
$ws.Range('E51').Value = '  +9.18%  '
